$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (No. 10) for the "Future work" part, to be written by Thao
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Future work"
$ws.Range("C11").Value = "Thao"

# Match the table formatting (font / alignment / borders) already used by the
# previous data row, without introducing any new cell styles
$ws.Range("A10:C10").Copy() | Out-Null
$ws.Range("A11:C11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the active selection to match the final state of the sheet
$ws.Range("M20").Select() | Out-Null
